$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet "GroupWithDifferentTradingDates" ---
# It is a copy of worksheet "2" (same BrokerageNotes, same group/style),
# placed right after worksheet "1", with the trading date of the middle
# row changed so the group now has different trading dates.
$source = $wb.Worksheets.Item("2")
$target = $wb.Worksheets.Item("1")
$source.Copy($null, $target)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "GroupWithDifferentTradingDates"

# Row 3 used to share the same trading date (39757) as rows 2 and 4.
# Change it to 39758 so the three BrokerageNotes in the group no longer
# all share the same TradingDate.
$newSheet.Range("A3").Value = 39758

# Reflect the cursor position left on the new sheet.
$newSheet.Range("A3").Select()

# --- 2. Update cursor position left on worksheet "1" ---
$ws1 = $wb.Worksheets.Item("1")
$ws1.Range("G2").Select()

# --- 3. Update cursor position left on worksheet "2" ---
$ws2 = $wb.Worksheets.Item("2")
$ws2.Range("A7").Select()

# --- 4. Make the new sheet the active / selected tab ---
$newSheet.Activate()
